# "Add files via upload" -- content update to 제품 백로그.xlsx
# 강한빛 is added as an additional worker on the "데이터 수집" backlog item (ID 2),
# and its estimated-work-time formula note is updated from 1*2*0.8=1.6 to 2*2*0.8=3.2,
# on both the product backlog sheet and the sprint plan sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("제품 백로그")
$ws2 = $wb.Worksheets.Item("스프린트 계획")

# --- Sheet "제품 백로그" : row for ID 2 ("데이터 수집") ---
$ws1.Range("G3").Value = "강한빛, 정혜미"
$ws1.Range("D3").Value = "2*2*0.8=3.2"

# Re-apply the existing center/center alignment to the edited worker cell
$ws1.Range("G3").HorizontalAlignment = -4108
$ws1.Range("G3").VerticalAlignment = -4108

# --- Sheet "스프린트 계획" : row for ID 2 ("데이터 수집") ---
$ws2.Range("C4").Value = "강한빛, 정혜미"
$ws2.Range("D4").Value = "2*2*0.8=3.2"

# Re-apply the existing center/center alignment to the edited worker cell
$ws2.Range("C4").HorizontalAlignment = -4108
$ws2.Range("C4").VerticalAlignment = -4108

# --- Restore/update the cursor selections left on each sheet before saving ---
$ws2.Activate()
$ws2.Range("G18").Select()

$ws1.Activate()
$ws1.Range("C18:C19").Select()

# --- Cosmetic window setting (tab-bar / horizontal scrollbar divider ratio) ---
$excel.ActiveWindow.TabRatio = 582
